# Refresh of the "cryptos" price/volume table (GitHub Actions scheduled update).
# All Price/Volume(1h) cells are stored as plain text in the sheet (e.g. "94.65",
# "  -4.40%  "), so for any Price cell whose new text would otherwise be parsed
# as a genuine number by Excel on entry (single decimal point, no thousands
# separator) we briefly force Text format before writing it, then clear that
# formatting again so the cell ends up with its original (default) style but
# keeps the literal text value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textForceRefs = @("D5","D6","D9","D10","D11","D13","D16","D19","D22","D23","D24","D28","D30","D31","D32","D34","D35","D36","D37","D38","D39","D41","D42","D46","D49","D50","D51")
foreach ($r in $textForceRefs) { $ws.Range($r).NumberFormat = "@" }

$ws.Range("D2").Value = "41.513.62"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "2.462.17"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "310.66"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "94.65"
$ws.Range("E6").Value = "  -4.40%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  -3.74%  "
$ws.Range("D10").Value = "33.69"
$ws.Range("E10").Value = "  -5.74%  "
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "6.96"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").Value = "2.841.93"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "2.469.39"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").Value = "14.53"
$ws.Range("E16").Value = "  -7.47%  "
$ws.Range("E17").Value = "  -3.91%  "
$ws.Range("D18").Value = "41.504.04"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  -6.22%  "
$ws.Range("D20").Value = "0.0₃0917"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("E21").Value = "  -4.91%  "
$ws.Range("D22").Value = "69.67"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").Value = "236.80"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("D24").Value = "2.77"
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -4.59%  "
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -4.90%  "
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("D30").Value = "36.32"
$ws.Range("E30").Value = "  -7.32%  "
$ws.Range("D31").Value = "153.84"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").Value = "5.60"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "2.55"
$ws.Range("E34").Value = "  -8.20%  "
$ws.Range("D35").Value = "0.0754"
$ws.Range("E35").Value = "  -4.63%  "
$ws.Range("D36").Value = "3.01"
$ws.Range("E36").Value = "  -4.53%  "
$ws.Range("D37").Value = "17.31"
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("D38").Value = "1.87"
$ws.Range("E38").Value = "  -6.90%  "
$ws.Range("D39").Value = "0.105"
$ws.Range("E39").Value = "  -5.34%  "
$ws.Range("E40").Value = "  -3.40%  "
$ws.Range("D41").Value = "4.01"
$ws.Range("E41").Value = "  -6.57%  "
$ws.Range("D42").Value = "21.33"
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D44").Value = "1.982.23"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("E45").Value = "  -4.15%  "
$ws.Range("D46").Value = "3.06"
$ws.Range("E46").Value = "  -7.07%  "
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").Value = "2.701.61"
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("D49").Value = "76.45"
$ws.Range("E49").Value = "  -5.52%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "69.57"
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("D51").Value = "97.57"
$ws.Range("E51").Value = "  -3.53%  "

# Restore default (general) cell formatting now that the literal text values
# are committed, so the only lasting change is the cell content.
foreach ($r in $textForceRefs) { $ws.Range($r).ClearFormats() }
